$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14) on the "Repayment schedule"
# sheet, shifting the existing "Late"/"Outstanding"/"Paid Date"/"Disbursement"
# columns one place to the right.
$wsSchedule.Columns.Item(14).Insert()

# Match the width of the newly inserted column to its left neighbour (column M)
# so it ends up with an explicit custom width of 11, rather than the sheet
# default width.
$wsSchedule.Columns.Item(14).ColumnWidth = $wsSchedule.Columns.Item(13).ColumnWidth

# Move the active sheet / selection from "NewLoanInput" to "Repayment schedule",
# and update the selected cell on that sheet.
$wsSchedule.Activate() | Out-Null
$wsSchedule.Range("R9").Select() | Out-Null
